# إضافة حدث جديد في Card14
#
# The "Card14" sheet is an append-only service log: row 14 was the last
# completed entry (13\8\2024 - تشحيم + صيانة) and still had several of its
# tracking columns (Min_Tones .. Revolving flats(o), Event) left blank.
# A new semi-annual-maintenance event (12\12\2024) needs to be logged as the
# next row, and row 14's blank tracking cells get backfilled with "nan" to
# match the convention used by every other completed row on the sheet.
#
# We insert a fresh blank row at position 14 (shifting the existing, already
# fully-populated row 14 down to row 15) rather than simply appending a new
# row 15 from scratch. That keeps row 15's still-blank tracking cells as
# real (empty-string) cells carried over from the original row, exactly as
# row 14 itself looked before this edit - a plain `.Value = ""` assignment
# on a brand-new cell would instead clear/remove it entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card14")

# Push the existing (complete) row 14 down to row 15, leaving a blank row 14
# in its place for the previously-existing event.
$ws.Cells.Item(14, 1).EntireRow.Insert()

# --- Re-create row 14 (the 13\8\2024 event) in the newly inserted row ---
# Column A ("card") is stored as text "14" throughout the sheet, so force a
# text number-format before assigning, then drop back to the default
# "Normal" style so no stray style index is left behind on the cell.
$a14 = $ws.Cells.Item(14, 1)
$a14.NumberFormat = "@"
$a14.Value = "14"
$a14.Style = "Normal"

$ws.Range("B14").Value = "nan"
$ws.Range("C14").Value = "nan"
$ws.Range("D14").Value = "nan"
$ws.Range("E14").Value = "nan"
$ws.Range("F14").Value = "nan"
$ws.Range("G14").Value = "nan"
$ws.Range("H14").Value = "nan"
$ws.Range("I14").Value = "nan"
$ws.Range("J14").Value = "nan"
$ws.Range("K14").Value = "nan"
$ws.Range("L14").Value = "13\8\2024"
$ws.Range("M14").Value = "nan"
$ws.Range("N14").Value = "تم تشحيم المكنه بالكامل +عمل صيانه"
$ws.Range("O14").Value = "تيم العمل"

# --- Row 15 (the shifted-down row) now becomes the new event ---
# Its "card" column (A15) already carries the text "14" from the shift, and
# the tracking columns B15:K15/M15 stay blank, exactly as required.
$ws.Range("L15").Value = "12\12\2024"
$ws.Range("N15").Value = "تم عمل صيانه نصف سنويه"
$ws.Range("O15").Value = "تيم العمل"
